# Regenerate sval data to filter save games.
# Updates the numeric B:G columns (rows 2-11) on the active sheet with
# the newly computed values. Column G is the row sum of B:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    3  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    5  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 4.429675500412797)
    6  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 17.45944343273191)
    7  = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 1, 10.08332054870323)
    8  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    9  = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 1.642425054193055)
    10 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    11 = @(0.6753301551942219, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 0, 13.27839780809431)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
